$d = $word.ActiveDocument

# --- Change 1: merge "Dallas" + ", TX" into a single run "Dallas, TX" ---
$d.Content.Find.Execute("Dallas, TX", $true, $false, $false, $false, $false, $true, 1, $false, "Dallas, TX", 2) | Out-Null

# --- Change 2: correct Dean's Honor Roll semester count 6 -> 4 ---
$rng = $d.Content
$rng.Find.Execute("(6 semesters)", $true, $false, $false, $false, $false, $true, 1, $false, "(4 semesters)", 2) | Out-Null

# Locate the "4" that was just inserted so the resulting run is split into
# three pieces ("(", "4", " semesters), 2021-2023") the way the source edit
# left it, instead of being re-coalesced into a single run.
$rng2 = $d.Content
$rng2.Find.Execute("(4 semesters)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$four = $d.Range($rng2.Start + 1, $rng2.Start + 2)
$four.Font.Italic = 1
$four.Font.Italic = 0
